# Weekly update: prepend two new price-report rows (Primera / Segunda)
# for "Apio" at row 395, pushing the existing rows 395:420 down to 397:422.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 395 (shifts old rows 395-420 down to 397-422,
# mirroring the existing row formatting the way Excel does on row insert).
$ws.Rows("395:396").Insert()

# --- New row 395 (Primera) ---
$ws.Cells.Item(395, 1).Value = 9
$ws.Cells.Item(395, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(395, 3).Value = "Metropolitana"
$ws.Cells.Item(395, 4).Value = 45013
$ws.Cells.Item(395, 5).Value = 13
$ws.Cells.Item(395, 6).Value = 100112017
$ws.Cells.Item(395, 7).Value = "Apio"
$ws.Cells.Item(395, 8).Value = "Americana (o)"
$ws.Cells.Item(395, 9).Value = "Primera"
$ws.Cells.Item(395, 10).Value = 70
$ws.Cells.Item(395, 11).Value = 7000
$ws.Cells.Item(395, 12).Value = 8000
$ws.Cells.Item(395, 13).Value = 7500
$ws.Cells.Item(395, 14).Value = "`$/docena de matas"
$ws.Cells.Item(395, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(395, 16).Value = 1250
$ws.Cells.Item(395, 17).Value = 6
$ws.Cells.Item(395, 18).Value = "Hortaliza"

# --- New row 396 (Segunda) ---
$ws.Cells.Item(396, 1).Value = 9
$ws.Cells.Item(396, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(396, 3).Value = "Metropolitana"
$ws.Cells.Item(396, 4).Value = 45013
$ws.Cells.Item(396, 5).Value = 13
$ws.Cells.Item(396, 6).Value = 100112017
$ws.Cells.Item(396, 7).Value = "Apio"
$ws.Cells.Item(396, 8).Value = "Americana (o)"
$ws.Cells.Item(396, 9).Value = "Segunda"
$ws.Cells.Item(396, 10).Value = 43
$ws.Cells.Item(396, 11).Value = 6000
$ws.Cells.Item(396, 12).Value = 6000
$ws.Cells.Item(396, 13).Value = 6000
$ws.Cells.Item(396, 14).Value = "`$/docena de matas"
$ws.Cells.Item(396, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(396, 16).Value = 1000
$ws.Cells.Item(396, 17).Value = 6
$ws.Cells.Item(396, 18).Value = "Hortaliza"
